$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the Price/Volume columns being updated so
# that numeric-looking strings (e.g. "1.004", "21.40") are preserved exactly
# instead of being auto-converted/normalized to numbers by Excel, then restore
# the original (unstyled) cell style so formatting is left untouched.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.430.01'
$ws.Range("E2").Value = '  -2.82%  '
$ws.Range("D3").Value = '1.742.57'
$ws.Range("E3").Value = '  -3.36%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '322.42'
$ws.Range("E5").Value = '  -4.35%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").Value = '0.4216'
$ws.Range("E7").Value = '  -9.76%  '
$ws.Range("D8").Value = '0.3578'
$ws.Range("E8").Value = '  -6.22%  '
$ws.Range("D9").Value = '45.46'
$ws.Range("E9").Value = '  +0.74%  '
$ws.Range("D10").Value = '0.07412'
$ws.Range("E10").Value = '  -2.56%  '
$ws.Range("D11").Value = '1.111'
$ws.Range("E11").Value = '  -3.29%  '
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.11%  '
$ws.Range("D13").Value = '21.40'
$ws.Range("E13").Value = '  -4.43%  '
$ws.Range("D14").Value = '6.088'
$ws.Range("E14").Value = '  -3.86%  '
$ws.Range("D15").Value = '7.174'
$ws.Range("E15").Value = '  -3.72%  '
$ws.Range("D16").Value = '1.742.74'
$ws.Range("E16").Value = '  -3.30%  '
$ws.Range("D17").Value = '0.00001064'
$ws.Range("E17").Value = '  -2.60%  '
$ws.Range("D18").Value = '87.21'
$ws.Range("E18").Value = '  +6.73%  '
$ws.Range("D19").Value = '0.06180'
$ws.Range("E19").Value = '  -8.09%  '
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").Value = '16.83'
$ws.Range("E21").Value = '  -3.20%  '
$ws.Range("D22").Value = '6.090'
$ws.Range("E22").Value = '  -4.95%  '
$ws.Range("D23").Value = '0.5253'
$ws.Range("E23").Value = '  -5.50%  '
$ws.Range("D24").Value = '27.489.86'
$ws.Range("E24").Value = '  -2.66%  '
$ws.Range("D25").Value = '11.60'
$ws.Range("E25").Value = '  -2.10%  '
$ws.Range("D26").Value = '2.325'
$ws.Range("E26").Value = '  -3.35%  '
$ws.Range("D27").Value = '20.43'
$ws.Range("E27").Value = '  -1.58%  '
$ws.Range("D28").Value = '151.84'
$ws.Range("E28").Value = '  -1.32%  '
$ws.Range("D29").Value = '2.359'
$ws.Range("E29").Value = '  -0.41%  '
$ws.Range("D30").Value = '1.939.97'
$ws.Range("E30").Value = '  -3.45%  '
$ws.Range("D31").Value = '126.01'
$ws.Range("E31").Value = '  -5.60%  '
$ws.Range("D32").Value = '1.203'
$ws.Range("E32").Value = '  -4.05%  '
$ws.Range("D33").Value = '5.665'
$ws.Range("E33").Value = '  -3.15%  '
$ws.Range("D34").Value = '0.09131'
$ws.Range("E34").Value = '  -5.01%  '
$ws.Range("D35").Value = '3.690'
$ws.Range("E35").Value = '  -8.58%  '
$ws.Range("D36").Value = '12.62'
$ws.Range("E36").Value = '  +4.46%  '
$ws.Range("D37").Value = '0.02283'
$ws.Range("E37").Value = '  -2.95%  '
$ws.Range("D38").Value = '5.083'
$ws.Range("E38").Value = '  -3.22%  '
$ws.Range("D39").Value = '0.2125'
$ws.Range("E39").Value = '  -4.68%  '
$ws.Range("D40").Value = '0.06079'
$ws.Range("E40").Value = '  -4.42%  '
$ws.Range("D41").Value = '0.6391'
$ws.Range("E41").Value = '  -3.50%  '
$ws.Range("D42").Value = '1.190'
$ws.Range("E42").Value = '  -3.65%  '
$ws.Range("D43").Value = '1.422'
$ws.Range("E43").Value = '  -5.21%  '
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").Value = '7.878'
$ws.Range("E45").Value = '  -4.26%  '
$ws.Range("D46").Value = '13.72'
$ws.Range("E46").Value = '  -3.79%  '
$ws.Range("D47").Value = '3.719'
$ws.Range("E47").Value = '  -2.98%  '
$ws.Range("D48").Value = '0.5856'
$ws.Range("E48").Value = '  -4.31%  '
$ws.Range("D49").Value = '124.72'
$ws.Range("E49").Value = '  -4.34%  '
$ws.Range("D50").Value = '1.948'
$ws.Range("E50").Value = '  -4.16%  '
$ws.Range("D51").Value = '0.06848'
$ws.Range("E51").Value = '  -4.37%  '

$ws.Range("D2:E51").Style = "Normal"
